$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column C (old C->D, D->E, E->F, F->G, G->H)
$ws.Columns("C").Insert()

# 2. Column widths (values chosen so the pixel-quantized COM result matches
# the target stored widths as closely as possible)
$ws.Columns("A").ColumnWidth = 22.33
$ws.Columns("C").ColumnWidth = 29.67
$ws.Columns("E").ColumnWidth = 17.67

# 3. New cell values - order matters for shared-string table layout
$ws.Range("D9").Value = "artikel 32bb"
$ws.Range("D10").Value = "artikel 3.92b"
$ws.Range("D11").Value = "artikel 32bc"
$ws.Range("B11").Value = "onderdeel b"
$ws.Range("F22").Value = "In hetderde lid (ACHTUNG: leerzeichen fehlt)"
$ws.Range("E8").Value = "punt twee "
$ws.Range("C1").Value = "subonderdeel"
$ws.Range("C5").Value = "subonderdeel 2"
$ws.Range("D12").Value = "van artikel 9a"
$ws.Range("B12").Value = "Onderdeel B"
$ws.Range("B13").Value = "Onderdeel b"
$ws.Range("A9").Value = "artikelen I tot en met III"
$ws.Range("C6").Value = "onder 2 (?) (59)"
$ws.Range("A18").Value = "Note: until amends[60]"

# 4. Repeat aanhef/considerans values on row 5
$ws.Range("G5").Value = $ws.Range("G1").Value2
$ws.Range("H5").Value = $ws.Range("H1").Value2

# 5. Selection
$ws.Range("B20").Select() | Out-Null
